$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.699.63"
$ws.Range("E2").Value = "  +3.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.081.56"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.03"
$ws.Range("E5").Value = "  +4.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.12"
$ws.Range("E6").Value = "  +7.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.074.65"
$ws.Range("E8").Value = "  +3.55%  "
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  +4.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("E13").Value = "  +5.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.22"
$ws.Range("E14").Value = "  +6.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.572.16"
$ws.Range("E15").Value = "  +2.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.719.22"
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.079.60"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  +5.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "487.42"
$ws.Range("E20").Value = "  +7.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.56"
$ws.Range("E21").Value = "  +4.17%  "
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.23"
$ws.Range("E23").Value = "  +6.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.44"
$ws.Range("E24").Value = "  +5.72%  "
$ws.Range("E25").Value = "  +5.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  +5.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  +6.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("E29").Value = "  +10.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.12"
$ws.Range("E31").Value = "  +3.39%  "
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.87"
$ws.Range("E33").Value = "  +11.33%  "
$ws.Range("E34").Value = "  +9.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.64"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.01"
$ws.Range("E36").Value = "  +4.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "469.49"
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.198.94"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0820"
$ws.Range("E39").Value = "  +6.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0399"
$ws.Range("E40").Value = "  +4.80%  "
$ws.Range("E41").Value = "  +5.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.25"
$ws.Range("E42").Value = "  +4.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  +7.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.87"
$ws.Range("E44").Value = "  +12.96%  "
$ws.Range("E45").Value = "  +4.74%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  +8.26%  "
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0524"
$ws.Range("E49").Value = "  +5.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.55"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("E51").Value = "  +7.55%  "
